$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 593, shifting existing rows 593:619 down to 594:620.
$ws.Rows.Item(593).Insert()

# Populate the newly inserted row 593 with the new weekly record.
$ws.Range("A593").Value = 3
$ws.Range("B593").Value = "Femacal de La Calera"
$ws.Range("C593").Value = "Coquimbo"
$ws.Range("D593").Value = 45147
$ws.Range("E593").Value = 5
$ws.Range("F593").Value = 100112040
$ws.Range("G593").Value = "Cilantro"
$ws.Range("H593").Value = "Sin especificar"
$ws.Range("I593").Value = "Primera"
$ws.Range("J593").Value = 50
$ws.Range("K593").Value = 4000
$ws.Range("L593").Value = 4000
$ws.Range("M593").Value = 4000
$ws.Range("N593").Value = "$/docena de atados (3 kilos)"
$ws.Range("O593").Value = "Provincia de Quillota"
$ws.Range("P593").Value = 1333
$ws.Range("Q593").Value = 3
$ws.Range("R593").Value = "Hortaliza"
